$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after the "little_saved" dialogue line (row 26) to add a
# modAffinity(+100) action, pushing everything below down by one row.
$ws.Rows(27).Insert()
$ws.Range("H27").Clear()
$ws.Rows(27).RowHeight = 91
$ws.Range("D27").Value = "modAffinity"
$ws.Range("E27").Value = 100

# Restore the MAX() formula reference which the row insert shifted out of range.
$ws.Range("H2").Formula = "=MAX(H4:H1048576)"

# Insert a new row after the "little_dead" second dialogue line (originally
# row 33, now row 34) to add a modAffinity(-200) action, pushing the trailing
# "end" marker down by one row.
$ws.Rows(35).Insert()
$ws.Range("H35").Clear()
$ws.Rows(35).RowHeight = 13.8
$ws.Range("D35").Value = "modAffinity"
$ws.Range("E35").Value = -200

# Restore the MAX() formula reference again after the second row insert.
$ws.Range("H2").Formula = "=MAX(H4:H1048576)"

# Leave the cursor on the new last action cell, matching where the editor
# ended up after making this change.
$ws.Range("H37").Select()
